$d = $word.ActiveDocument

function Find-ParagraphIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq $text) {
            return $i
        }
    }
    return -1
}

function Replace-ParagraphText($oldText, $innerXml) {
    $idx = Find-ParagraphIndex $oldText
    $p = $d.Paragraphs($idx)
    $rng = $p.Range
    $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + $innerXml + "</w:p>"
    $rng.InsertXML($xml) | Out-Null
}

# Heading1 title run (only the first run changes; the second run
# "Gameplay and Features" that immediately follows it is untouched)
Replace-ParagraphText "Play Glory of Egypt for Free - Exciting Slot GameGameplay and Features" "<w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:t>Play Glory of Egypt Free - A Captivating Ancient Egypt Slot Game</w:t></w:r><w:r><w:t>Gameplay and Features</w:t></w:r>"

# "What we like" bullet list
Replace-ParagraphText "Above-average RTP of 96.05%" "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Well-designed graphics and animations</w:t></w:r>"
Replace-ParagraphText "Captivating graphics and animations" "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Captivating ancient Egypt theme</w:t></w:r>"
Replace-ParagraphText "Wilds, Scatters, and free spins offer many chances to win" "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Variety of winning features including Wilds, Scatters, and free spins</w:t></w:r>"
Replace-ParagraphText "Bonus gambling game offers an opportunity to double winnings" "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Chance to double winnings with bonus gambling game</w:t></w:r>"

# "What we don't like" bullet list
Replace-ParagraphText "Limited 10 ways to win" "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Limited number of paylines (10 ways to win)</w:t></w:r>"
Replace-ParagraphText "Lack of bonus features" "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>High volatility may not appeal to all players</w:t></w:r>"

# Bold title recap near the end of the document
Replace-ParagraphText "Play Glory of Egypt for Free - Exciting Slot Game" "<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Glory of Egypt Free - A Captivating Ancient Egypt Slot Game</w:t></w:r>"

# Italic meta description at the very end of the document
Replace-ParagraphText "Read our review of Glory of Egypt, a captivating slot game set in ancient Egypt. Play this game for free and enjoy Wilds, Scatters, and free spins." "<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Play Glory of Egypt free and experience the captivating world of ancient Egypt with exciting features.</w:t></w:r>"
